$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift product_price/deliver_date/no_of_items/special_note/payment_method
# one column to the right (L..P -> M..Q) and insert a new package_id value
# in L2, working right-to-left so source cells aren't clobbered before
# they're read.
$ws.Range("Q2").Value = 1
$ws.Range("P2").Value = "delivery safely plz"
$ws.Range("O2").Value = 1
$ws.Range("N2").NumberFormat = "m/d/yy"
$ws.Range("N2").Value = 44310

# recipient_number / recipient_number_two become text-formatted phone
# numbers (keep the leading zero).
$ws.Range("B2:C2").NumberFormat = "@"
$ws.Range("B2").Value = "01917720168"
$ws.Range("C2").Value = "01921307315"

$ws.Range("M2").NumberFormat = "0"
$ws.Range("M2").Value = 500
$ws.Range("L2").Value = 1

# Column width tweaks: columns O (no_of_items) and P (special_note) get
# genuinely wider, no-longer-autofit columns. (The other columns' widths
# shift by only thousandths in the source diff -- an artifact of the file
# being re-saved by a newer Excel build's font metrics -- so they're left
# alone here.)
$ws.Columns("O").ColumnWidth = 12
$ws.Columns("P").ColumnWidth = 22

# View changes: scroll so column J is the left-most visible column and
# select P8.
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("P8").Select()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
